$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------

# Grab the run layout (empty run + bold run) from the paragraph that
# currently holds the bold "Play Cash Stampede..." text at the end of
# the document, so the freshly inserted paragraph gets an identical
# <w:r/><w:r><w:rPr><w:b/></w:rPr>... run structure.
$paraCount = $d.Paragraphs.Count
$boldTitlePara = $d.Paragraphs($paraCount - 1)
$boldTitlePara.Range.Copy()

$titlePara = $d.Paragraphs.First
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$metaRange.Collapse(1)            # wdCollapseStart
$metaRange.Paste()

# $metaRange now sits at the (collapsed) insertion point; the pasted
# paragraph content starts right there. Re-point a working range at
# the pasted bold text ("Play Cash Stampede Free Today | Affordable
# Betting and Wild Wins" -> 64 characters) and swap its text for the
# "Meta description" label.
$labelRange = $d.Range($metaRange.Start, $metaRange.Start + 64)
$labelRange.Text = "Meta description"

# Append the (non-bold) description text right after the label run.
$labelRange.Collapse(0)           # wdCollapseEnd
$labelRange.InsertAfter(": Try your luck with Cash Stampede, a safari-themed slot game with affordable betting options, bonus features, and Wild Cash Stampede wins. Play for free now.")

# ------------------------------------------------------------------
# 2) Drop the now-duplicated bold "Play Cash Stampede..." paragraph
#    near the end of the document.
# ------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($paraCount - 1)
$dupTitlePara.Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the italic "Try your luck..." blurb at the very end with
#    the image-generation prompt text (formatting/run stays italic).
#    Scope the Find/Replace to just this last paragraph's Range so it
#    can't also match the (textually similar) meta-description blurb
#    that now lives near the top of the document.
# ------------------------------------------------------------------
$oldText = "Try your luck with Cash Stampede, a safari-themed slot game with affordable betting options, bonus features, and Wild Cash Stampede wins. Play for free now."
$newText = "Prompt: Create a feature image for Cash Stampede in a cartoon style. The image should feature a happy Maya warrior with glasses. The warrior should be depicted holding a lasso and riding on a majestic elephant with the other African animals in the background. The background should be a beautiful sunset with silhouettes of trees and grass. The overall image should give off a sense of adventure and excitement. The dimensions for the image should be 1080 x 1080 pixels."
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
